$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1422834409152216, 0.1945146085939768, 0.1112616825839581, 3.552228279157816, 0.002512203253484145, 0.6729235202556652, 0.3118672535452305)
    3 = @(0.1330520610950003, 0.1834900375207127, 0.1020511602801264, 3.339602839758271, 0.002518234767026999, 0.6195835765612685, 0.2868217909492969)
    4 = @(0.127461414131659, 0.1766980199772092, 0.09644059510097946, 3.209681841570074, 0.002522125131106571, 0.5877094959811302, 0.2717181039332317)
    5 = @(0.1252026710787959, 0.1739239692793859, 0.09416515595338382, 3.156889713247068, 0.002523757690054559, 0.5749380251172056, 0.2656311202340262)
    6 = @(0.124828786180359, 0.1734629485269465, 0.09378797029230412, 3.148132577307138, 0.002524031631692878, 0.5728303798741763, 0.2646244501584079)
    7 = @(0.127430872994907, 0.1766606341098509, 0.09640986399422502, 3.208969263507385, 0.002522146956954078, 0.5875363788548498, 0.2716357392808035)
    8 = @(0.1390843891877438, 0.1907177859272053, 0.1080764038067414, 3.478779487155407, 0.002514244220676529, 0.654348065628426, 0.3031739716825612)
    9 = @(0.1625523457950493, 0.2181256561564453, 0.1313258163943587, 4.01324113747205, 0.002500222236873449, 0.7924689014371609, 0.3672502588964122)
    10 = @(0.1801727575816301, 0.2381999421613159, 0.1486581019323623, 4.409707060345312, 0.002490807792589794, 0.8984926994328362, 0.4157661074989107)
    11 = @(0.188271806009439, 0.2473262583074813, 0.1566029697515248, 4.59101171408588, 0.002486715087166504, 0.9477613465671482, 0.4381675421705893)
    12 = @(0.1913507336319356, 0.2507819380697356, 0.159620575943336, 4.659812645498675, 0.002485192407743596, 0.9665711620894797, 0.4466994094565706)
    13 = @(0.190687098092738, 0.2500376983520027, 0.1589702714787933, 4.644988563044762, 0.002485519139669545, 0.9625132756019639, 0.4448597222477133)
    14 = @(0.1885248708795899, 0.2476105611937953, 0.1568510462948609, 4.596669059108137, 0.002486589272397973, 0.9493057571929171, 0.4388684759351946)
    15 = @(0.1872020061416748, 0.2461238522488429, 0.1555541506560658, 4.567091068171635, 0.002487248288574109, 0.9412357843961274, 0.4352050753321777)
    16 = @(0.1796451442374973, 0.2376034500542801, 0.1481401325771472, 4.397878190222315, 0.002491079067945555, 0.8952940684449686, 0.4143088975326208)
    17 = @(0.1750306202038843, 0.232375415740421, 0.1436075933603504, 4.294320716541762, 0.002493477656161113, 0.8673786789345286, 0.4015755370759493)
    18 = @(0.1723843213926273, 0.2293678308000864, 0.1410062514205066, 4.234846149805009, 0.002494875152766903, 0.8514199561226121, 0.3942827529565776)
    19 = @(0.1714896790711435, 0.2283494019405765, 0.1401264413903291, 4.214724115491407, 0.002495351399247807, 0.846033239094595, 0.3918188456845826)
    20 = @(0.1755210314088771, 0.2329320027376411, 0.1440895019021156, 4.305335321576081, 0.002493220472116939, 0.8703402061092618, 0.4029277962492799)
    21 = @(0.189159643938595, 0.2483234728487105, 0.1574732652123458, 4.610857673461567, 0.002486274213472738, 0.9531809481510152, 0.4406269118802442)
    22 = @(0.1981431743126336, 0.2583814412078311, 0.1662733025450223, 4.811382040668263, 0.002481892537235903, 1.008214806957966, 0.4655512816510452)
    23 = @(0.1933420991336163, 0.2530132468488375, 0.1615715817565544, 4.704278107440018, 0.002484216712623751, 0.9787593127746277, 0.4522220922796123)
    24 = @(0.1752992957816986, 0.2326803757721052, 0.1438716171398369, 4.300355426123872, 0.002493336687368956, 0.8690010189721136, 0.4023163534003942)
    25 = @(0.1561374202045442, 0.2107245250022345, 0.1249939602571217, 3.868026814523347, 0.002503858844182484, 0.7543210310070947, 0.3496690614117952)
}

$cols = @("B","D","E","F","G","K","M")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

$wb.Save()